$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to snake_case machine-readable names ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Normalize municipality/state names: capitalize linking words (de, del, la, las, los, el, y) ---
$ws.Range('B7').Value = 'Pabellón De Arteaga'
$ws.Range('B8').Value = 'Rincón De Romos'
$ws.Range('B9').Value = 'San Francisco De Los Romo'
$ws.Range('B30').Value = 'Amatenango De La Frontera'
$ws.Range('B35').Value = 'Chiapa De Corzo'
$ws.Range('B53').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B56').Value = 'San Cristóbal De Las Casas'
$ws.Range('B85').Value = 'Guadalupe Y Calvo'
$ws.Range('B87').Value = 'Hidalgo Del Parral'
$ws.Range('B101').Value = 'San Francisco De Borja'
$ws.Range('B102').Value = 'San Francisco Del Oro'
$ws.Range('B110').Value = 'Valle De Zaragoza'
$ws.Range('B132').Value = 'Villa De Álvarez'
$ws.Range('A134').Value = 'Ciudad De México'
$ws.Range('B149').Value = 'Coneto De Comonfort'
$ws.Range('B160').Value = 'Nombre De Dios'
$ws.Range('B163').Value = 'Pánuco De Coronado'
$ws.Range('B170').Value = 'San Luis Del Cordero'
$ws.Range('A177').Value = 'Estado De México'
$ws.Range('B177').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B180').Value = 'Almoloya De Alquisiras'
$ws.Range('B181').Value = 'Almoloya De Juárez'
$ws.Range('B186').Value = 'Atizapán De Zaragoza'
$ws.Range('B190').Value = 'Chapa De Mota'
$ws.Range('B193').Value = 'Coacalco De Berriozábal'
$ws.Range('B198').Value = 'Ecatepec De Morelos'
$ws.Range('B202').Value = 'Ixtapan De La Sal'
$ws.Range('B216').Value = 'Naucalpan De Juárez'
$ws.Range('B225').Value = 'San Felipe Del Progreso'
$ws.Range('B226').Value = 'San Martín De Las Pirámides'
$ws.Range('B227').Value = 'San Simón De Guerrero'
$ws.Range('B237').Value = 'Tenango Del Valle'
$ws.Range('B247').Value = 'Tlalnepantla De Baz'
$ws.Range('B253').Value = 'Valle De Bravo'
$ws.Range('B254').Value = 'Villa De Allende'
$ws.Range('B255').Value = 'Villa Del Carbón'
$ws.Range('B266').Value = 'San Miguel De Allende'
$ws.Range('B267').Value = 'Apaseo El Alto'
$ws.Range('B268').Value = 'Apaseo El Grande'
$ws.Range('B275').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B285').Value = 'Purísima Del Rincón'
$ws.Range('B289').Value = 'San Diego De La Unión'
$ws.Range('B291').Value = 'San Francisco Del Rincón'
$ws.Range('B293').Value = 'San Luis De La Paz'
$ws.Range('B294').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B296').Value = 'Silao De La Victoria'
$ws.Range('B299').Value = 'Valle De Santiago'
$ws.Range('B303').Value = 'Acapulco De Juárez'
$ws.Range('B305').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B306').Value = 'Alcozauca De Guerrero'
$ws.Range('B310').Value = 'Atenango Del Río'
$ws.Range('B312').Value = 'Atoyac De Álvarez'
$ws.Range('B313').Value = 'Ayutla De Los Libres'
$ws.Range('B316').Value = 'Buenavista De Cuéllar'
$ws.Range('B317').Value = 'Chilapa De Álvarez'
$ws.Range('B318').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B319').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B324').Value = 'Coyuca De Benítez'
$ws.Range('B325').Value = 'Coyuca De Catalán'
$ws.Range('B328').Value = 'Cuetzala Del Progreso'
$ws.Range('B329').Value = 'Cutzamala De Pinzón'
$ws.Range('B335').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B336').Value = 'Iguala De La Independencia'
$ws.Range('B338').Value = 'Zihuatanejo De Azueta'
$ws.Range('B340').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B343').Value = 'Mártir De Cuilapan'
$ws.Range('B356').Value = 'Taxco De Alarcón'
$ws.Range('B358').Value = 'Técpan De Galeana'
$ws.Range('B360').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B362').Value = 'Tixtla De Guerrero'
$ws.Range('B365').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B366').Value = 'Tlapa De Comonfort'
$ws.Range('B375').Value = 'Agua Blanca De Iturbide'
$ws.Range('B379').Value = 'Atotonilco De Tula'
$ws.Range('B380').Value = 'Atotonilco El Grande'
$ws.Range('B383').Value = 'Cuautepec De Hinojosa'
$ws.Range('B387').Value = 'Huasca De Ocampo'
$ws.Range('B389').Value = 'Huejutla De Reyes'
$ws.Range('B391').Value = 'Jacala De Ledezma'
$ws.Range('B395').Value = 'Mineral Del Monte'
$ws.Range('B396').Value = 'Mixquiahuala De Juárez'
$ws.Range('B398').Value = 'Pachuca De Soto'
$ws.Range('B401').Value = 'Progreso De Obregón'
$ws.Range('B404').Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range('B407').Value = 'Tenango De Doria'
$ws.Range('B409').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B410').Value = 'Tezontepec De Aldama'
$ws.Range('B416').Value = 'Tula De Allende'
$ws.Range('B417').Value = 'Tulancingo De Bravo'
$ws.Range('B418').Value = 'Villa De Tezontepec'
$ws.Range('B420').Value = 'Zacualtipán De Ángeles'
$ws.Range('B425').Value = 'Acatlán De Juárez'
$ws.Range('B426').Value = 'Ahualulco De Mercado'
$ws.Range('B430').Value = 'Atemajac De Brizuela'
$ws.Range('B432').Value = 'Atotonilco El Alto'
$ws.Range('B433').Value = 'Autlán De Navarro'
$ws.Range('B439').Value = 'Cañadas De Obregón'
$ws.Range('B443').Value = 'Concepción De Buenos Aires'
$ws.Range('B444').Value = 'Cuautitlán De García Barragán'
$ws.Range('B450').Value = 'Encarnación De Díaz'
$ws.Range('B457').Value = 'Huejuquilla El Alto'
$ws.Range('B458').Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range('B459').Value = 'Ixtlahuacán Del Río'
$ws.Range('B468').Value = 'Lagos De Moreno'
$ws.Range('B475').Value = 'Ojuelos De Jalisco'
$ws.Range('B480').Value = 'San Diego De Alejandría'
$ws.Range('B481').Value = 'San Juan De Los Lagos'
$ws.Range('B482').Value = 'San Martín De Bolaños'
$ws.Range('B484').Value = 'San Miguel El Alto'
$ws.Range('B485').Value = 'San Sebastián Del Oeste'
$ws.Range('B488').Value = 'Talpa De Allende'
$ws.Range('B489').Value = 'Tamazula De Gordiano'
$ws.Range('B495').Value = 'Teocuitatlán De Corona'
$ws.Range('B496').Value = 'Tepatitlán De Morelos'
$ws.Range('B499').Value = 'Tizapán El Alto'
$ws.Range('B500').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B509').Value = 'Unión De San Antonio'
$ws.Range('B510').Value = 'Unión De Tula'
$ws.Range('B515').Value = 'Yahualica De González Gallo'
$ws.Range('B516').Value = 'Zacoalco De Torres'
$ws.Range('B519').Value = 'Zapotitlán De Vadillo'
$ws.Range('B520').Value = 'Zapotlán Del Rey'
$ws.Range('B521').Value = 'Zapotlán El Grande'
$ws.Range('B544').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B605').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B629').Value = 'Coatlán Del Río'
$ws.Range('B636').Value = 'Jonacatepec De Leandro Valle'
$ws.Range('B639').Value = 'Puente De Ixtla'
$ws.Range('B643').Value = 'Tlaltizapán De Zapata'
$ws.Range('B653').Value = 'Amatlán De Cañas'
$ws.Range('B654').Value = 'Bahía De Banderas'
$ws.Range('B657').Value = 'Ixtlán Del Río'
$ws.Range('B663').Value = 'Santa María Del Oro'
$ws.Range('B672').Value = 'Mier Y Noriega'
$ws.Range('B674').Value = 'San Nicolás De Los Garza'
$ws.Range('B676').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B679').Value = 'Ayoquezco De Aldama'
$ws.Range('B681').Value = 'Chalcatongo De Hidalgo'
$ws.Range('B683').Value = 'Coicoyán De Las Flores'
$ws.Range('B685').Value = 'Fresnillo De Trujano'
$ws.Range('B686').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B687').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B688').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B689').Value = 'Huautla De Jiménez'
$ws.Range('B690').Value = 'Ixtlán De Juárez'
$ws.Range('B691').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B697').Value = 'Mariscala De Juárez'
$ws.Range('B698').Value = 'Mártires De Tacubaya'
$ws.Range('B700').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B701').Value = 'Nejapa De Madero'
$ws.Range('B702').Value = 'Oaxaca De Juárez'
$ws.Range('B703').Value = 'Ocotlán De Morelos'
$ws.Range('B704').Value = 'Putla Villa De Guerrero'
$ws.Range('B715').Value = 'San Antonio De La Cal'
$ws.Range('B727').Value = 'San Juan Bautista Lo De Soto'
$ws.Range('B752').Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Range('B755').Value = 'Santa Ana Del Valle'
$ws.Range('B759').Value = 'Santa Inés Del Monte'
$ws.Range('B767').Value = 'Santa María Jalapa Del Marqués'
$ws.Range('B785').Value = 'Santo Domingo De Morelos'
$ws.Range('B790').Value = 'Sitio De Xitlapehua'
$ws.Range('B791').Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range('B792').Value = 'Teotitlán Del Valle'
$ws.Range('B793').Value = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
$ws.Range('B794').Value = 'Tlacolula De Matamoros'
$ws.Range('B795').Value = 'Villa De Etla'
$ws.Range('B796').Value = 'Villa De Tututepec'
$ws.Range('B797').Value = 'Villa De Zaachila'
$ws.Range('B798').Value = 'Villa Sola De Vega'
$ws.Range('B800').Value = 'Zimatlán De Álvarez'
$ws.Range('B827').Value = 'Cuayuca De Andrade'
$ws.Range('B837').Value = 'Huehuetlán El Grande'
$ws.Range('B840').Value = 'Ixcamilpa De Guerrero'
$ws.Range('B843').Value = 'Izúcar De Matamoros'
$ws.Range('B848').Value = 'Los Reyes De Juárez'
$ws.Range('B855').Value = 'Palmar De Bravo'
$ws.Range('B867').Value = 'San Nicolás De Los Ranchos'
$ws.Range('B872').Value = 'Tecali De Herrera'
$ws.Range('B879').Value = 'Tepanco De López'
$ws.Range('B883').Value = 'Tepexi De Rodríguez'
$ws.Range('B884').Value = 'Tetela De Ocampo'
$ws.Range('B889').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B895').Value = 'Totoltepec De Guerrero'
$ws.Range('B900').Value = 'Xayacatlán De Bravo'
$ws.Range('B906').Value = 'Amealco De Bonfil'
$ws.Range('B908').Value = 'Cadereyta De Montes'
$ws.Range('B912').Value = 'Pinal De Amoles'
$ws.Range('B915').Value = 'San Juan Del Río'
$ws.Range('B925').Value = 'Ciudad Del Maíz'
$ws.Range('B931').Value = 'Mexquitic De Carmona'
$ws.Range('B936').Value = 'San Ciro De Acosta'
$ws.Range('B938').Value = 'Santa María Del Río'
$ws.Range('B940').Value = 'Soledad De Graciano Sánchez'
$ws.Range('B945').Value = 'Villa De Arriaga'
$ws.Range('B990').Value = 'Nacozari De García'
$ws.Range('B1022').Value = 'Acuamanala De Miguel Hidalgo'
$ws.Range('B1024').Value = 'Amaxac De Guerrero'
$ws.Range('B1031').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B1032').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B1035').Value = 'Papalotla De Xicohténcatl'
$ws.Range('B1036').Value = 'San Pablo Del Monte'
$ws.Range('B1038').Value = 'Tepetitla De Lardizábal'
$ws.Range('B1041').Value = 'Tetla De La Solidaridad'
$ws.Range('B1056').Value = 'Amatlán De Los Reyes'
$ws.Range('B1060').Value = 'Boca Del Río'
$ws.Range('B1063').Value = 'Castillo De Teayo'
$ws.Range('B1073').Value = 'Cosamaloapan De Carpio'
$ws.Range('B1086').Value = 'Hueyapan De Ocampo'
$ws.Range('B1087').Value = 'Ignacio De La Llave'
$ws.Range('B1089').Value = 'Ixhuatlán De Madero'
$ws.Range('B1090').Value = 'Ixhuatlán Del Café'
$ws.Range('B1098').Value = 'Juchique De Ferrer'
$ws.Range('B1104').Value = 'Martínez De La Torre'
$ws.Range('B1106').Value = 'Medellín De Bravo'
$ws.Range('B1117').Value = 'Paso Del Macho'
$ws.Range('B1119').Value = 'Poza Rica De Hidalgo'
$ws.Range('B1126').Value = 'Sayula De Alemán'
$ws.Range('B1127').Value = 'Soledad De Doblado'
$ws.Range('B1147').Value = 'Vega De Alatorre'
$ws.Range('B1154').Value = 'Zozocolco De Hidalgo'
$ws.Range('B1178').Value = 'Cañitas De Felipe Pescador'
$ws.Range('B1180').Value = 'Concepción Del Oro'
$ws.Range('B1181').Value = 'El Plateado De Joaquín Amaro'
$ws.Range('B1193').Value = 'Mezquital Del Oro'
$ws.Range('B1196').Value = 'Moyahua De Estrada'
$ws.Range('B1197').Value = 'Nochistlán De Mejía'
$ws.Range('B1198').Value = 'Noria De Ángeles'
$ws.Range('B1208').Value = 'Teúl De González Ortega'
$ws.Range('B1209').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B1212').Value = 'Villa De Cos'

# --- Recalculated percentage values (float precision adjustments) ---
$ws.Range('D4').Value = 0.0009317296290022024
$ws.Range('D18').Value = 0.0009317296290022024
$ws.Range('D137').Value = 0.0009317296290022024
$ws.Range('D150').Value = 0.0009317296290022024
$ws.Range('D224').Value = 0.0009317296290022024
$ws.Range('D265').Value = 0.0009317296290022024
$ws.Range('D291').Value = 0.0009317296290022024
$ws.Range('D292').Value = 0.0009317296290022024
$ws.Range('D322').Value = 0.0009317296290022024
$ws.Range('D338').Value = 0.0009317296290022024
$ws.Range('D430').Value = 0.0009317296290022024
$ws.Range('D547').Value = 0.0009317296290022024
$ws.Range('D560').Value = 0.0009317296290022024
$ws.Range('D640').Value = 0.0009317296290022024
$ws.Range('D651').Value = 0.0009317296290022024
$ws.Range('D794').Value = 0.0009317296290022024
$ws.Range('D1016').Value = 0.0009317296290022024
$ws.Range('D1033').Value = 0.0009317296290022024
$ws.Range('D1053').Value = 0.0009317296290022024
$ws.Range('D1070').Value = 0.0009317296290022024
$ws.Range('D1206').Value = 0.0009317296290022024
$ws.Range('D134').Value = 0.009571404370658988
$ws.Range('D352').Value = 0.009147890902930712

# --- Remove trailing metadata rows (1219-1224) and shrink used range to A1:D1218 ---
$ws.Rows('1219:1224').Delete()

Write-Output ('Final dimension: ' + $ws.UsedRange.Address())
